# "finestra incidenza 7gg centrata su ultimo g"
# Recompute the 7-day moving sum (column C) and the 7-day moving incidence
# per 100,000 inhabitants (column D) so that the window of 7 days ends on
# (is centred on) the current/last day, i.e. C(r) = SUM(B(r-6):B(r)).
# Previously the window ended 3 days earlier (C(r) = SUM(B(r-9):B(r-3))).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$lastDataRow = 184
$windowSize = 7
$perInhabitants = 100000
$population = 2080

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $windowStart = $r - $windowSize + 1

    if ($windowStart -lt $firstDataRow) {
        # Not enough prior days to fill a full 7 day window -> leave blank
        $ws.Cells.Item($r, 3).ClearContents()
        $ws.Cells.Item($r, 4).ClearContents()
    }
    else {
        $sum = 0
        for ($i = $windowStart; $i -le $r; $i++) {
            $sum += $ws.Cells.Item($i, 2).Value2
        }
        $ws.Cells.Item($r, 3).Value2 = $sum
        $ws.Cells.Item($r, 4).Value2 = $sum * $perInhabitants / $population
    }
}

$wb.Save()
